# Auto-generated edit script: refreshes market-price / profit columns (H-N)
# on the Leve profit sheets, per the scheduled-runner data update.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3862.628
$ws.Range("J40").Value = 4203.875
$ws.Range("L40").Value = 4203.875
$ws.Range("N40").Value = -4553.875

$ws.Range("H51").Value = 42618
$ws.Range("I51").Value = 150000
$ws.Range("J51").Value = 31879.8
$ws.Range("K51").Value = 150000
$ws.Range("L51").Value = 31879.8
$ws.Range("M51").Value = -149516
$ws.Range("N51").Value = -32847.8

$ws.Range("H62").Value = 4301.3
$ws.Range("I62").Value = 3572.5715
$ws.Range("K62").Value = 3572.5715
$ws.Range("M62").Value = -2948.5715

$ws.Range("H65").Value = 4301.3
$ws.Range("I65").Value = 3572.5715
$ws.Range("K65").Value = 17862.8575
$ws.Range("M65").Value = -14742.8575

$ws.Range("H74").Value = 6782.8335
$ws.Range("I74").Value = 6782.8335
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 6782.8335
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -5846.8335
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 6782.8335
$ws.Range("I77").Value = 6782.8335
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 33914.1675
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -29234.1675
$ws.Range("N77").ClearContents()

$ws.Range("H86").Value = 4254.364
$ws.Range("I86").Value = 4929.933
$ws.Range("J86").Value = 2806.7144
$ws.Range("K86").Value = 4929.933
$ws.Range("L86").Value = 2806.7144
$ws.Range("M86").Value = -3806.933
$ws.Range("N86").Value = -5052.7144

$ws.Range("H89").Value = 4254.364
$ws.Range("I89").Value = 4929.933
$ws.Range("J89").Value = 2806.7144
$ws.Range("K89").Value = 24649.665
$ws.Range("L89").Value = 14033.572
$ws.Range("M89").Value = -19033.665
$ws.Range("N89").Value = -25265.572

$ws.Range("H106").Value = 2500
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()

$ws.Range("H138").Value = 2649.2092
$ws.Range("I138").Value = 643.8
$ws.Range("J138").Value = 3256.9092
$ws.Range("K138").Value = 1931.4
$ws.Range("L138").Value = 9770.7276
$ws.Range("M138").Value = 3208.6
$ws.Range("N138").Value = -20050.7276

$ws.Range("H141").Value = 3183.8262
$ws.Range("I141").Value = 3230.3809
$ws.Range("K141").Value = 9691.1427
$ws.Range("M141").Value = -4511.1427


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1910.7
$ws.Range("I2").Value = 1944
$ws.Range("K2").Value = 1944
$ws.Range("M2").Value = -1831

$ws.Range("H8").Value = 35861.2
$ws.Range("I8").Value = 30111.75
$ws.Range("J8").Value = 42432
$ws.Range("K8").Value = 30111.75
$ws.Range("L8").Value = 42432
$ws.Range("M8").Value = -29967.75
$ws.Range("N8").Value = -42720

$ws.Range("H34").Value = 500000
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H37").Value = 44500
$ws.Range("J37").Value = 44500
$ws.Range("L37").Value = 44500
$ws.Range("N37").Value = -45046

$ws.Range("H88").Value = 1305.1
$ws.Range("J88").Value = 1414.0714
$ws.Range("L88").Value = 1414.0714
$ws.Range("N88").Value = -2226.0714

$ws.Range("H91").Value = 1305.1
$ws.Range("J91").Value = 1414.0714
$ws.Range("L91").Value = 1414.0714
$ws.Range("N91").Value = -4222.0714

$ws.Range("H116").Value = 1910.7
$ws.Range("I116").Value = 1944
$ws.Range("K116").Value = 1944
$ws.Range("M116").Value = 350


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1910.7
$ws.Range("I3").Value = 1944
$ws.Range("K3").Value = 1944
$ws.Range("M3").Value = -1830

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 248.1
$ws.Range("I7").Value = 230.14285
$ws.Range("J7").Value = 290
$ws.Range("K7").Value = 230.14285
$ws.Range("L7").Value = 290
$ws.Range("M7").Value = -117.14285
$ws.Range("N7").Value = -516

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H110").Value = 49994.5
$ws.Range("J110").Value = 49994.5
$ws.Range("L110").Value = 49994.5
$ws.Range("N110").Value = -58174.5


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 951.63635
$ws.Range("I12").Value = 1859
$ws.Range("J12").Value = 611.375
$ws.Range("K12").Value = 5577
$ws.Range("L12").Value = 1834.125
$ws.Range("M12").Value = -5404
$ws.Range("N12").Value = -2180.125

$ws.Range("H37").Value = 85000
$ws.Range("J37").Value = 85000
$ws.Range("L37").Value = 255000
$ws.Range("N37").Value = -255224

$ws.Range("H38").Value = 51.75
$ws.Range("I38").Value = 15.666667
$ws.Range("J38").Value = 58.117645
$ws.Range("K38").Value = 47.000001
$ws.Range("L38").Value = 174.352935
$ws.Range("M38").Value = 299.999999
$ws.Range("N38").Value = -868.352935


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 262.30768
$ws.Range("I2").Value = 20.75
$ws.Range("K2").Value = 20.75
$ws.Range("M2").Value = 92.25

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1955.6875
$ws.Range("I16").Value = 1530.5385
$ws.Range("K16").Value = 1530.5385
$ws.Range("M16").Value = -1360.5385

$ws.Range("H46").Value = 3107.4285
$ws.Range("I46").Value = 2608.3157
$ws.Range("J46").Value = 4161.1113
$ws.Range("K46").Value = 2608.3157
$ws.Range("L46").Value = 4161.1113
$ws.Range("M46").Value = -2420.3157
$ws.Range("N46").Value = -4537.1113


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4276.3335
$ws.Range("I122").Value = 2455.9473
$ws.Range("J122").Value = 8599.75
$ws.Range("K122").Value = 7367.841899999999
$ws.Range("L122").Value = 25799.25
$ws.Range("M122").Value = -4917.841899999999
$ws.Range("N122").Value = -30699.25

